$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Columns F..I were relabelled (rotated): SessionLabel/Subject/Session/Dataset
# -> Subject/Session/Dataset/SessionLabel
$ws.Range("F1").Value = "Subject"
$ws.Range("G1").Value = "Session"
$ws.Range("H1").Value = "Dataset"
$ws.Range("I1").Value = "SessionLabel"

# Columns J..O gain a "DICOM:" prefix. Set in this exact order so the new
# shared-string entries land in the same order as the target file.
$ws.Range("J1").Value = "DICOM:Manufacturer"
$ws.Range("K1").Value = "DICOM:ManufacturerModelName"
$ws.Range("L1").Value = "DICOM:Modality"
$ws.Range("M1").Value = "DICOM:StationName"
$ws.Range("N1").Value = "DICOM:StudyDate"
$ws.Range("O1").Value = "DICOM:StudyDescription"

# --- Data rows ---
# For each "Scans" row that had Subject/Session/Dataset values in G/H/I,
# those values shift left into F/G/H (matching the new header order), and I
# is cleared.
$rows = @(
    @{ Row = 3;  F = "002304";  G = "20200312";  H = "Head_CT" },
    @{ Row = 4;  F = "002304";  G = "20200312";  H = "Head_CT" },
    @{ Row = 5;  F = "002304";  G = "20200312";  H = "Neck_CT" },
    @{ Row = 6;  F = "002304";  G = "20200312";  H = "Neck_CT" },
    @{ Row = 7;  F = "002304";  G = "20200312";  H = "Neck_CT" },
    @{ Row = 10; F = "397829";  G = "20190115";  H = "SomeCT" },
    @{ Row = 12; F = "397829";  G = "20200623";  H = "SomeCT" },
    @{ Row = 14; F = "397829";  G = "20210414";  H = "SomeCT" },
    @{ Row = 16; F = "038945";  G = "20200303";  H = "X-Rays" }
)

foreach ($r in $rows) {
    $ws.Range("F$($r.Row)").Value = $r.F
    $ws.Range("G$($r.Row)").Value = $r.G
    $ws.Range("H$($r.Row)").Value = $r.H
    $ws.Range("I$($r.Row)").ClearContents()
}

# --- Selection moved from M1 to O2 ---
$ws.Range("O2").Select() | Out-Null
